# Finish Enhanced Logout Reminder Logic on Shift Change document
#
# - F9  ("Not update attributes of job for UCO of equipment when ending
#         order") gains the "Synchronize Equipment Counter with Line
#         During End Order" document hyperlink (moved from F17).
# - F13 ("Disable button 'Current SKU' for AU detail") is marked "No Needed".
# - F17 ("Enhanced logout remind logic on shift change") now points to its
#        own finished document, "Enhanced Logout Reminder Logic on Shift
#        Change".
# - The active sheet view scrolls back up and the selection moves to C11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F9: add the "Synchronize Equipment Counter..." hyperlink (previously on F17).
# Copy F17's cell format (the bordered Hyperlink style) first, then set the formula.
$ws.Range("F17").Copy()
$ws.Range("F9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F9").Formula = '=HYPERLINK("./Synchronize%20Equipment%20Counter%20with%20Line%20During%20End%20Job.docx", "Synchronize Equipment Counter with Line During End Order")'

# F13: mark as "No Needed"
$ws.Range("F13").Value = "No Needed"

# F17: replace with the new, finished "Enhanced Logout Reminder Logic on Shift Change" hyperlink
$ws.Range("F17").Formula = '=HYPERLINK("./Enhanced%20Logout%20Reminder%20Logic%20on%20Shift%20Change.docx", "Enhanced Logout Reminder Logic on Shift Change")'

# Restore the selection as saved in the edited workbook (the sheet also
# scrolls so A7 is the top-left visible cell, but that pure view/scroll
# state isn't part of the persisted cell/formula data this script controls)
$ws.Activate()
$ws.Range("C11").Select()
